$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 new values
$ws.Range("A2").Value = 112181650
$ws.Range("B2").Value = 78699
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("Q2").Value = 756202
$ws.Range("R2").Value = 7291065

# Row 3 new values
$ws.Range("A3").Value = 112181583
$ws.Range("B3").Value = 89780
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = "Fläckporing"
$ws.Range("G3").Value = "Anthoporia albobrunnea"
$ws.Range("H3").Value = "(Romell) Karasiński & Niemelä"
$ws.Range("Q3").Value = 756188
$ws.Range("R3").Value = 7291007

# Row 4 change
$ws.Range("B4").Value = 77388
